# Apply the crypto price / 1h-volume-change refresh described by the commit
# ("Updated cryptos list ... with GitHub Actions"). Only the D (Price) and
# E (Volume(1h)) columns change, row by row, on the single data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Plain decimal-looking Price values (e.g. "209.56") would otherwise be
# auto-coerced to numbers by Excel on assignment; the source data keeps them
# as literal text (some prices use "." as a thousands separator too, e.g.
# "25.906.78", which already is not numeric-looking and needs no special
# handling). A leading apostrophe forces the new value to stay text, just
# like typing it into the cell would -- the apostrophe itself is not stored,
# only the digits/dots after it become the cell text.

$ws.Range("D2").Value = '25.906.78'
$ws.Range("E2").Value = '  -1.06%  '
$ws.Range("D3").Value = '1.634.22'
$ws.Range("E3").Value = '  -2.58%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = "'209.56"
$ws.Range("E5").Value = '  -0.82%  '
$ws.Range("D6").Value = "'0.5201"
$ws.Range("E6").Value = '  -0.81%  '
$ws.Range("D7").Value = "'1.002"
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = "'0.2559"
$ws.Range("E8").Value = '  -3.66%  '
$ws.Range("D9").Value = "'0.06235"
$ws.Range("E9").Value = '  -0.94%  '
$ws.Range("D10").Value = "'20.25"
$ws.Range("E10").Value = '  -5.49%  '
$ws.Range("D11").Value = "'0.07554"
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("D12").Value = '1.628.65'
$ws.Range("E12").Value = '  -2.86%  '
$ws.Range("D13").Value = "'4.352"
$ws.Range("E13").Value = '  -2.31%  '
$ws.Range("D14").Value = '1.858.40'
$ws.Range("E14").Value = '  -2.53%  '
$ws.Range("D15").Value = "'0.5409"
$ws.Range("E15").Value = '  -4.16%  '
$ws.Range("D16").Value = '0.0₅7900'
$ws.Range("E16").Value = '  -1.66%  '
$ws.Range("D17").Value = "'64.40"
$ws.Range("E17").Value = '  -3.34%  '
$ws.Range("D18").Value = '25.918.26'
$ws.Range("E18").Value = '  -1.29%  '
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = "'4.620"
$ws.Range("E20").Value = '  -4.44%  '
$ws.Range("D21").Value = "'184.00"
$ws.Range("D22").Value = "'9.987"
$ws.Range("E22").Value = '  -4.25%  '
$ws.Range("D23").Value = "'6.073"
$ws.Range("E23").Value = '  -2.16%  '
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = '  -0.10%  '
$ws.Range("D25").Value = "'145.39"
$ws.Range("E25").Value = '  -2.68%  '
$ws.Range("D26").Value = "'0.1198"
$ws.Range("E26").Value = '  -3.91%  '
$ws.Range("D27").Value = "'7.317"
$ws.Range("E27").Value = '  -3.62%  '
$ws.Range("D28").Value = "'15.46"
$ws.Range("E28").Value = '  -3.41%  '
$ws.Range("E29").Value = '  +1.03%  '
$ws.Range("D30").Value = "'0.05947"
$ws.Range("E30").Value = '  -3.73%  '
$ws.Range("E31").Value = '  -3.29%  '
$ws.Range("D32").Value = "'3.344"
$ws.Range("E32").Value = '  -2.77%  '
$ws.Range("D33").Value = "'3.334"
$ws.Range("E33").Value = '  -4.62%  '
$ws.Range("D34").Value = "'1.603"
$ws.Range("E34").Value = '  -1.96%  '
$ws.Range("D35").Value = "'0.9680"
$ws.Range("E35").Value = '  -3.42%  '
$ws.Range("D36").Value = "'2.382"
$ws.Range("E36").Value = '  -1.00%  '
$ws.Range("D37").Value = "'2.727"
$ws.Range("E37").Value = '  -0.60%  '
$ws.Range("E38").Value = '  -4.26%  '
$ws.Range("D39").Value = "'0.01590"
$ws.Range("E39").Value = '  -1.40%  '
$ws.Range("E40").Value = '  -0.45%  '
$ws.Range("D41").Value = "'0.8385"
$ws.Range("E41").Value = '  -3.57%  '
$ws.Range("D42").Value = '1.025.32'
$ws.Range("E42").Value = '  -5.24%  '
$ws.Range("D43").Value = "'5.649"
$ws.Range("E43").Value = '  -7.37%  '
$ws.Range("D44").Value = "'99.57"
$ws.Range("D45").Value = '1.783.29'
$ws.Range("D46").Value = '0.0₈105'
$ws.Range("E46").Value = '  -5.48%  '
$ws.Range("D47").Value = "'0.9994"
$ws.Range("E47").Value = '  +0.13%  '
$ws.Range("E48").Value = '  -3.82%  '
$ws.Range("D49").Value = "'7.976"
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("D50").Value = "'0.05179"
$ws.Range("E50").Value = '  -1.17%  '
$ws.Range("D51").Value = "'0.4228"
$ws.Range("E51").Value = '  -0.71%  '
